$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "24.798.30"
$ws.Range("E2").Value2 = "  -4.31%  "

$ws.Range("D3").Value2 = "1.647.05"
$ws.Range("E3").Value2 = "  -6.91%  "

$ws.Range("D4").Value2 = "'1.005"
$ws.Range("E4").Value2 = "  +0.47%  "

$ws.Range("D5").Value2 = "'308.47"
$ws.Range("E5").Value2 = "  -2.79%  "

$ws.Range("D6").Value2 = "'1.004"
$ws.Range("E6").Value2 = "  +0.82%  "

$ws.Range("D7").Value2 = "'0.3637"
$ws.Range("E7").Value2 = "  -5.54%  "

$ws.Range("B8").Value2 = "Cardano"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value2 = "'0.3295"
$ws.Range("E8").Value2 = "  -9.90%  "

$ws.Range("B9").Value2 = "OKB"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value2 = "'46.80"
$ws.Range("E9").Value2 = "  -8.63%  "

$ws.Range("D10").Value2 = "'1.131"
$ws.Range("E10").Value2 = "  -8.87%  "

$ws.Range("D11").Value2 = "'0.07067"
$ws.Range("E11").Value2 = "  -8.28%  "

$ws.Range("D12").Value2 = "'1.005"
$ws.Range("E12").Value2 = "  +0.84%  "

$ws.Range("D13").Value2 = "'6.052"
$ws.Range("E13").Value2 = "  -7.29%  "

$ws.Range("D14").Value2 = "'19.70"
$ws.Range("E14").Value2 = "  -10.27%  "

$ws.Range("D15").Value2 = "'6.666"
$ws.Range("E15").Value2 = "  -6.85%  "

$ws.Range("D16").Value2 = "1.642.06"
$ws.Range("E16").Value2 = "  -6.94%  "

$ws.Range("D17").Value2 = "'0.00001065"
$ws.Range("E17").Value2 = "  -9.06%  "

$ws.Range("D18").Value2 = "'0.06591"
$ws.Range("E18").Value2 = "  -4.00%  "

$ws.Range("E19").Value2 = "  +0.79%  "

$ws.Range("D20").Value2 = "'79.25"
$ws.Range("E20").Value2 = "  -9.64%  "

$ws.Range("D21").Value2 = "'16.29"
$ws.Range("E21").Value2 = "  -8.44%  "

$ws.Range("D22").Value2 = "'6.013"
$ws.Range("E22").Value2 = "  -8.42%  "

$ws.Range("D23").Value2 = "'12.19"
$ws.Range("E23").Value2 = "  -5.11%  "

$ws.Range("D24").Value2 = "24.681.11"
$ws.Range("E24").Value2 = "  -4.54%  "

$ws.Range("D25").Value2 = "'2.419"
$ws.Range("E25").Value2 = "  -0.49%  "

$ws.Range("D26").Value2 = "'2.525"
$ws.Range("E26").Value2 = "  -16.03%  "

$ws.Range("D27").Value2 = "'148.36"
$ws.Range("E27").Value2 = "  -4.55%  "

$ws.Range("D28").Value2 = "'19.21"
$ws.Range("E28").Value2 = "  -7.90%  "

$ws.Range("D29").Value2 = "'127.79"
$ws.Range("E29").Value2 = "  -5.68%  "

$ws.Range("D30").Value2 = "1.826.09"
$ws.Range("E30").Value2 = "  -6.91%  "

$ws.Range("D31").Value2 = "'1.094"
$ws.Range("E31").Value2 = "  -9.12%  "

$ws.Range("D32").Value2 = "'4.119"
$ws.Range("E32").Value2 = "  -4.45%  "

$ws.Range("D33").Value2 = "'6.064"
$ws.Range("E33").Value2 = "  -16.69%  "

$ws.Range("D34").Value2 = "'1.727"
$ws.Range("E34").Value2 = "  -4.97%  "

$ws.Range("D35").Value2 = "'0.08444"
$ws.Range("E35").Value2 = "  -3.47%  "

$ws.Range("D36").Value2 = "'12.59"
$ws.Range("E36").Value2 = "  -10.88%  "

$ws.Range("D37").Value2 = "'5.214"
$ws.Range("E37").Value2 = "  -8.78%  "

$ws.Range("D38").Value2 = "'0.06193"
$ws.Range("E38").Value2 = "  -8.85%  "

$ws.Range("D39").Value2 = "'0.02285"
$ws.Range("E39").Value2 = "  -8.75%  "

$ws.Range("D40").Value2 = "'0.2090"
$ws.Range("E40").Value2 = "  -7.12%  "

$ws.Range("D41").Value2 = "'1.217"
$ws.Range("E41").Value2 = "  -6.71%  "

$ws.Range("D42").Value2 = "'8.291"
$ws.Range("E42").Value2 = "  -12.06%  "

$ws.Range("D43").Value2 = "'0.6062"
$ws.Range("E43").Value2 = "  -8.49%  "

$ws.Range("D44").Value2 = "'1.004"
$ws.Range("E44").Value2 = "  +0.77%  "

$ws.Range("D45").Value2 = "'3.753"
$ws.Range("E45").Value2 = "  -4.42%  "

$ws.Range("D46").Value2 = "'13.04"
$ws.Range("E46").Value2 = "  -8.99%  "

$ws.Range("D47").Value2 = "'0.5760"
$ws.Range("E47").Value2 = "  -10.07%  "

$ws.Range("D48").Value2 = "'123.67"
$ws.Range("E48").Value2 = "  -7.79%  "

$ws.Range("D49").Value2 = "'1.976"
$ws.Range("E49").Value2 = "  -9.57%  "

$ws.Range("D50").Value2 = "'0.07057"
$ws.Range("E50").Value2 = "  -6.09%  "

$ws.Range("D51").Value2 = "'74.95"
$ws.Range("E51").Value2 = "  -7.83%  "
